# Update the "State" column (B) for a set of rows from "yes" to "best".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsToUpdate = @(2, 3, 4, 10, 16, 23, 26, 27, 36)
foreach ($row in $rowsToUpdate) {
    $ws.Cells.Item($row, 2).Value = "best"
}

# Move the selection/active cell to C16 (and clear any frozen/scrolled top-left cell).
$ws.Range("C16").Select()
